$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.811.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.16%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.814.79"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.55%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.50%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.54"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.91%  "

# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.30%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.38%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "39.66"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +9.48%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.292"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.01%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0674"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.18%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1000"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.52%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.075.38"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.49%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.807.15"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.21%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.05"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.04%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.637"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.48%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.729.38"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.94%  "

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.25%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.40"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.02%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.59"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.58%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0770"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.25%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.19"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.56%  "

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.41%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.12"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.41%  "

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.76%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.55"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.02%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.77"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.74%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.59"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.30%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.09%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.49%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.61%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.77"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.51%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.89"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.61%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0515"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.73%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.83"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.97%  "

# Row 35
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.319.01"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.61%  "

# Row 36
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.644"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.62%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.14%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.38"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.80%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.23%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.24"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.10%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "83.06"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.87%  "

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.40"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.50%  "

# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.44"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.58%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.951"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.04%  "

# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.80"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.63%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0519"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.81%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.976.72"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.52%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.75"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.72%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.41%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.22"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.66%  "

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.45%  "
